$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the columns beyond W (X1:AQ19) that are no longer used ---
$ws.Range("X1:AQ19").Clear()

# --- Row 2: header HKL labels for columns C:W ---
$ws.Range("C2").Value = "[4, 0, 0]"
$ws.Range("D2").Value = "[4, 2, 0]"
$ws.Range("E2").Value = "[2, 2, 0]"
$ws.Range("F2").Value = "[2, 0, 0]"
$ws.Range("G2").Value = "[3, 3, 3]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[1, 1, 1]"
$ws.Range("J2").Value = "[3, 1, 1]"
$ws.Range("K2").Value = "[3, 3, 1]"
$ws.Range("L2").Value = "[4, 2, 2]"
$ws.Range("M2").Value = "[5, 1, 1]"
$ws.Range("N2").Value = "1Pair-A"
$ws.Range("O2").Value = "1Pair-B"
$ws.Range("P2").Value = "2Pairs-A"
$ws.Range("Q2").Value = "2Pairs-B"
$ws.Range("R2").Value = "3Pairs-A"
$ws.Range("S2").Value = "3Pairs-B"
$ws.Range("T2").Value = "3Pairs-C"
$ws.Range("U2").Value = "4Pairs"
$ws.Range("V2").Value = "5A4F"
$ws.Range("W2").Value = "MaxUnique"

# --- Rows 3-15: existing scheme rows, values re-permuted per new HKL column order ---
# Row 3: BT8Hex_2.5
$ws.Range("B3").Value = "BT8Hex_2.5"
$ws.Range("C3").Value = 0.994578759144797
$ws.Range("D3").Value = 0.9985135374786613
$ws.Range("E3").Value = 1.000735844379216
$ws.Range("F3").Value = 0.994578759144797
$ws.Range("G3").Value = 1.004759161331947
$ws.Range("H3").Value = 1.004759161331947
$ws.Range("I3").Value = 1.004759161331947
$ws.Range("J3").Value = 0.9987826297686138
$ws.Range("K3").Value = 1.001951357723442
$ws.Range("L3").Value = 1.001695887902029
$ws.Range("M3").Value = 0.9963595945147065
$ws.Range("N3").Value = 1.004759161331947
$ws.Range("O3").Value = 1.000735844379216
$ws.Range("P3").Value = 0.9976573017620065
$ws.Range("Q3").Value = 0.9997592370739149
$ws.Range("R3").Value = 1.00002458828532
$ws.Range("S3").Value = 0.9980324110975424
$ws.Range("T3").Value = 1.00002458828532
$ws.Range("U3").Value = 0.9997140986561432
$ws.Range("V3").Value = 1.000723111191304
$ws.Range("W3").Value = 0.9996720965304264

# Row 4: BT8Hex_5
$ws.Range("B4").Value = "BT8Hex_5"
$ws.Range("C4").Value = 0.9895132465174489
$ws.Range("D4").Value = 0.9971386508381694
$ws.Range("E4").Value = 1.001429042347536
$ws.Range("F4").Value = 0.9895132465174489
$ws.Range("G4").Value = 1.009165273922513
$ws.Range("H4").Value = 1.009165273922513
$ws.Range("I4").Value = 1.009165273922513
$ws.Range("J4").Value = 0.9976561951306122
$ws.Range("K4").Value = 1.003768280793965
$ws.Range("L4").Value = 1.003274744147699
$ws.Range("M4").Value = 0.9929683814368394
$ws.Range("N4").Value = 1.009165273922513
$ws.Range("O4").Value = 1.001429042347536
$ws.Range("P4").Value = 0.9954711444324924
$ws.Range("Q4").Value = 0.9995426187390741
$ws.Range("R4").Value = 1.000035854262499
$ws.Range("S4").Value = 0.996199494665199
$ws.Range("T4").Value = 1.000035854262499
$ws.Range("U4").Value = 0.9994409394795276
$ws.Range("V4").Value = 1.001385806368125
$ws.Range("W4").Value = 0.9993642268918479

# Row 5: BT8Hex_10
$ws.Range("B5").Value = "BT8Hex_10"
$ws.Range("C5").Value = 0.9798574642977054
$ws.Range("D5").Value = 0.9945559213465666
$ws.Range("E5").Value = 1.002849307318035
$ws.Range("F5").Value = 0.9798574642977054
$ws.Range("G5").Value = 1.017459532888245
$ws.Range("H5").Value = 1.017459532888245
$ws.Range("I5").Value = 1.017459532888245
$ws.Range("J5").Value = 0.9954901231685662
$ws.Range("K5").Value = 1.007258923844104
$ws.Range("L5").Value = 1.006236316749687
$ws.Range("M5").Value = 0.9864989958150467
$ws.Range("N5").Value = 1.017459532888245
$ws.Range("O5").Value = 1.002849307318035
$ws.Range("P5").Value = 0.9913533858078701
$ws.Range("Q5").Value = 0.9991697152433006
$ws.Range("R5").Value = 1.000055434834662
$ws.Range("S5").Value = 0.9927322982614356
$ws.Range("T5").Value = 1.000055434834662
$ws.Range("U5").Value = 0.9989141069181378
$ws.Range("V5").Value = 1.002623192112159
$ws.Range("W5").Value = 0.9987758231784943

# Row 6: BT8Hex_15
$ws.Range("B6").Value = "BT8Hex_15"
$ws.Range("C6").Value = 0.9704314596833274
$ws.Range("D6").Value = 0.9920464006424506
$ws.Range("E6").Value = 1.004249859488042
$ws.Range("F6").Value = 0.9704314596833274
$ws.Range("G6").Value = 1.025567568935174
$ws.Range("H6").Value = 1.025567568935174
$ws.Range("I6").Value = 1.025567568935174
$ws.Range("J6").Value = 0.9933638219349538
$ws.Range("K6").Value = 1.010666503664677
$ws.Range("L6").Value = 1.00911059658296
$ws.Range("M6").Value = 0.9801833677198277
$ws.Range("N6").Value = 1.025567568935174
$ws.Range("O6").Value = 1.004249859488042
$ws.Range("P6").Value = 0.9873406595856846
$ws.Range("Q6").Value = 0.9988068407114978
$ws.Range("R6").Value = 1.000082962702181
$ws.Range("S6").Value = 0.9893483803687744
$ws.Range("T6").Value = 1.000082962702181
$ws.Range("U6").Value = 0.9984031775103742
$ws.Range("V6").Value = 1.003836055795334
$ws.Range("W6").Value = 0.9982024473314265

# Row 7: Spiral2.5
$ws.Range("B7").Value = "Spiral2.5"
$ws.Range("C7").Value = 0.9994816263125663
$ws.Range("D7").Value = 0.9999165894080103
$ws.Range("E7").Value = 1.000163140017173
$ws.Range("F7").Value = 0.9994816263125663
$ws.Range("G7").Value = 1.000286493689811
$ws.Range("H7").Value = 1.000286493689811
$ws.Range("I7").Value = 1.000286493689811
$ws.Range("J7").Value = 0.9998904371583699
$ws.Range("K7").Value = 1.000196720399643
$ws.Range("L7").Value = 1.00011182456786
$ws.Range("M7").Value = 0.9996680232426394
$ws.Range("N7").Value = 1.000286493689811
$ws.Range("O7").Value = 1.000163140017173
$ws.Range("P7").Value = 0.9998223831648695
$ws.Range("Q7").Value = 1.000026788587771
$ws.Range("R7").Value = 0.9999770866731833
$ws.Range("S7").Value = 0.9998450678293697
$ws.Range("T7").Value = 0.9999770866731833
$ws.Range("U7").Value = 0.9999554242944799
$ws.Range("V7").Value = 1.000021638173546
$ws.Range("W7").Value = 0.9999643568495091

# Row 8: Spiral5
$ws.Range("B8").Value = "Spiral5"
$ws.Range("C8").Value = 0.9985199703430819
$ws.Range("D8").Value = 0.9997309736284997
$ws.Range("E8").Value = 1.00041559969905
$ws.Range("F8").Value = 0.9985199703430819
$ws.Range("G8").Value = 1.00090598303344
$ws.Range("H8").Value = 1.00090598303344
$ws.Range("I8").Value = 1.00090598303344
$ws.Range("J8").Value = 0.9996845941242937
$ws.Range("K8").Value = 1.000555583492638
$ws.Range("L8").Value = 1.000346382739877
$ws.Range("M8").Value = 0.9990443035715149
$ws.Range("N8").Value = 1.00090598303344
$ws.Range("O8").Value = 1.00041559969905
$ws.Range("P8").Value = 0.9994677850210663
$ws.Range("Q8").Value = 1.000050096911672
$ws.Range("R8").Value = 0.9999471843585243
$ws.Range("S8").Value = 0.999540054722142
$ws.Range("T8").Value = 0.9999471843585243
$ws.Range("U8").Value = 0.9998815367999666
$ws.Range("V8").Value = 1.000086426046661
$ws.Range("W8").Value = 0.9999004238290494

# Row 9: Spiral7.5
$ws.Range("B9").Value = "Spiral7.5"
$ws.Range("C9").Value = 0.9978459984854007
$ws.Range("D9").Value = 0.999556110967893
$ws.Range("E9").Value = 1.000521102935452
$ws.Range("F9").Value = 0.9978459984854007
$ws.Range("G9").Value = 1.00146595307938
$ws.Range("H9").Value = 1.00146595307938
$ws.Range("I9").Value = 1.00146595307938
$ws.Range("J9").Value = 0.9995364499552731
$ws.Range("K9").Value = 1.000798974765334
$ws.Range("L9").Value = 1.000549921187226
$ws.Range("M9").Value = 0.9985954666067383
$ws.Range("N9").Value = 1.00146595307938
$ws.Range("O9").Value = 1.000521102935452
$ws.Range("P9").Value = 0.9991835507104265
$ws.Range("Q9").Value = 1.000028776445363
$ws.Range("R9").Value = 0.9999443515000777
$ws.Range("S9").Value = 0.999301183792042
$ws.Range("T9").Value = 0.9999443515000777
$ws.Range("U9").Value = 0.9998423761138765
$ws.Range("V9").Value = 1.000167091506977
$ws.Range("W9").Value = 0.9998587472478373

# Row 10: Spiral10
$ws.Range("B10").Value = "Spiral10"
$ws.Range("C10").Value = 0.9954224894740417
$ws.Range("D10").Value = 0.9991144816031706
$ws.Range("E10").Value = 1.001192175443521
$ws.Range("F10").Value = 0.9954224894740417
$ws.Range("G10").Value = 1.002953482658046
$ws.Range("H10").Value = 1.002953482658046
$ws.Range("I10").Value = 1.002953482658046
$ws.Range("J10").Value = 0.9990238812647662
$ws.Range("K10").Value = 1.001704435994626
$ws.Range("L10").Value = 1.001120792656758
$ws.Range("M10").Value = 0.9970331297453559
$ws.Range("N10").Value = 1.002953482658046
$ws.Range("O10").Value = 1.001192175443521
$ws.Range("P10").Value = 0.9983073324587812
$ws.Range("Q10").Value = 1.000108028354143
$ws.Range("R10").Value = 0.9998560491918695
$ws.Range("S10").Value = 0.9985461820607763
$ws.Range("T10").Value = 0.9998560491918695
$ws.Range("U10").Value = 0.9996480072100937
$ws.Range("V10").Value = 1.000309102299684
$ws.Range("W10").Value = 0.9996956086050357

# Row 11: Spiral15
$ws.Range("B11").Value = "Spiral15"
$ws.Range("C11").Value = 0.9922161249489228
$ws.Range("D11").Value = 0.9983137330700799
$ws.Range("E11").Value = 1.001744250046807
$ws.Range("F11").Value = 0.9922161249489228
$ws.Range("G11").Value = 1.005530231307421
$ws.Range("H11").Value = 1.005530231307421
$ws.Range("I11").Value = 1.005530231307421
$ws.Range("J11").Value = 0.9983214527582972
$ws.Range("K11").Value = 1.002868241759702
$ws.Range("L11").Value = 1.002061656645142
$ws.Range("M11").Value = 0.9949057595459186
$ws.Range("N11").Value = 1.005530231307421
$ws.Range("O11").Value = 1.001744250046807
$ws.Range("P11").Value = 0.9969801874978647
$ws.Range("Q11").Value = 1.000032851402552
$ws.Range("R11").Value = 0.9998302021010502
$ws.Range("S11").Value = 0.9974272759180088
$ws.Range("T11").Value = 0.9998302021010502
$ws.Range("U11").Value = 0.999453014765362
$ws.Range("V11").Value = 1.000668458073774
$ws.Range("W11").Value = 0.9994951812602864

# Row 12: OffsetF45
$ws.Range("B12").Value = "OffsetF45"
$ws.Range("C12").Value = 1.114359926117516
$ws.Range("D12").Value = 1.000461040254294
$ws.Range("E12").Value = 0.9355091172145319
$ws.Range("F12").Value = 1.114359926117516
$ws.Range("G12").Value = 0.9860814526514601
$ws.Range("H12").Value = 0.9860814526514601
$ws.Range("I12").Value = 0.9860814526514601
$ws.Range("J12").Value = 1.022939026056447
$ws.Range("K12").Value = 0.9534045963821534
$ws.Range("L12").Value = 0.9911820440698969
$ws.Range("M12").Value = 1.068629186633957
$ws.Range("N12").Value = 0.9860814526514601
$ws.Range("O12").Value = 0.9355091172145319
$ws.Range("P12").Value = 1.024934521666024
$ws.Range("Q12").Value = 0.9792240716354892
$ws.Range("R12").Value = 1.011983498661169
$ws.Range("S12").Value = 1.024269356462831
$ws.Range("T12").Value = 1.011983498661169
$ws.Range("U12").Value = 1.014722380509989
$ws.Range("V12").Value = 1.008994194938283
$ws.Range("W12").Value = 1.009070798672532

# Row 13: OffsetA45
$ws.Range("B13").Value = "OffsetA45"
$ws.Range("C13").Value = 1.013023478209936
$ws.Range("D13").Value = 1.01186609296541
$ws.Range("E13").Value = 1.012259343999107
$ws.Range("F13").Value = 1.013023478209936
$ws.Range("G13").Value = 0.9655899923054061
$ws.Range("H13").Value = 0.9655899923054061
$ws.Range("I13").Value = 0.9655899923054061
$ws.Range("J13").Value = 1.003084459309294
$ws.Range("K13").Value = 0.9972410691309167
$ws.Range("L13").Value = 0.9883139550583141
$ws.Range("M13").Value = 1.010568862437565
$ws.Range("N13").Value = 0.9655899923054061
$ws.Range("O13").Value = 1.012259343999107
$ws.Range("P13").Value = 1.012641411104521
$ws.Range("Q13").Value = 1.0076719016542
$ws.Range("R13").Value = 0.9969576048381494
$ws.Range("S13").Value = 1.009455760506112
$ws.Range("T13").Value = 0.9969576048381494
$ws.Range("U13").Value = 0.9984893184559355
$ws.Range("V13").Value = 0.9919094532258296
$ws.Range("W13").Value = 1.000243406676993

# Row 14: OffsetFTD
$ws.Range("B14").Value = "OffsetFTD"
$ws.Range("C14").Value = 0.9495556907948551
$ws.Range("D14").Value = 1.029713845490989
$ws.Range("E14").Value = 1.073895747723929
$ws.Range("F14").Value = 0.9495556907948551
$ws.Range("G14").Value = 0.9232440436940611
$ws.Range("H14").Value = 0.9232440436940611
$ws.Range("I14").Value = 0.9232440436940611
$ws.Range("J14").Value = 0.9933444082733895
$ws.Range("K14").Value = 1.024771104951278
$ws.Range("L14").Value = 0.9783665009512694
$ws.Range("M14").Value = 0.9783002134435266
$ws.Range("N14").Value = 0.9232440436940611
$ws.Range("O14").Value = 1.073895747723929
$ws.Range("P14").Value = 1.011725719259392
$ws.Range("Q14").Value = 1.033620077998659
$ws.Range("R14").Value = 0.9822318274042816
$ws.Range("S14").Value = 1.005598615597391
$ws.Range("T14").Value = 0.9822318274042816
$ws.Range("U14").Value = 0.9850099726215586
$ws.Range("V14").Value = 0.9726567868360592
$ws.Range("W14").Value = 0.9938989444154123

# Row 15: OffsetATD
$ws.Range("B15").Value = "OffsetATD"
$ws.Range("C15").Value = 1.017742000111958
$ws.Range("D15").Value = 0.9869753866395347
$ws.Range("E15").Value = 0.9708447351522942
$ws.Range("F15").Value = 1.017742000111958
$ws.Range("G15").Value = 1.033861511597765
$ws.Range("H15").Value = 1.033861511597765
$ws.Range("I15").Value = 1.033861511597765
$ws.Range("J15").Value = 1.00171894190294
$ws.Range("K15").Value = 0.9913194689752158
$ws.Range("L15").Value = 1.009575558954179
$ws.Range("M15").Value = 1.006636380876238
$ws.Range("N15").Value = 1.033861511597765
$ws.Range("O15").Value = 0.9708447351522942
$ws.Range("P15").Value = 0.9942933676321261
$ws.Range("Q15").Value = 0.986281838527617
$ws.Range("R15").Value = 1.007482748954006
$ws.Range("S15").Value = 0.9967685590557306
$ws.Range("T15").Value = 1.007482748954006
$ws.Range("U15").Value = 1.006041797191239
$ws.Range("V15").Value = 1.011605740072544
$ws.Range("W15").Value = 1.002334248026266

# --- Rows 16-19: NEW Holden2.5/5/10/15 rows (A column needs value + style copy) ---
# Row 16: Holden2.5
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("C16").Value = 0.8930939820618167
$ws.Range("D16").Value = 0.9717067039018414
$ws.Range("E16").Value = 1.016026947532025
$ws.Range("F16").Value = 0.8930939820618167
$ws.Range("G16").Value = 1.090913023433971
$ws.Range("H16").Value = 1.090913023433971
$ws.Range("I16").Value = 1.090913023433971
$ws.Range("J16").Value = 0.9761689459463058
$ws.Range("K16").Value = 1.038605957713575
$ws.Range("L16").Value = 1.032606032654188
$ws.Range("M16").Value = 0.9285293965360057
$ws.Range("N16").Value = 1.090913023433971
$ws.Range("O16").Value = 1.016026947532025
$ws.Range("P16").Value = 0.9545604647969209
$ws.Range("Q16").Value = 0.9960979467391655
$ws.Range("R16").Value = 1.000011317675938
$ws.Range("S16").Value = 0.9617632918467159
$ws.Range("T16").Value = 1.000011317675938
$ws.Range("U16").Value = 0.9940507247435296
$ws.Range("V16").Value = 1.013423184481618
$ws.Range("W16").Value = 0.9934563737224662

# Row 17: Holden5
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Holden5"
$ws.Range("C17").Value = 0.9099865788991316
$ws.Range("D17").Value = 0.9784218905287612
$ws.Range("E17").Value = 1.017051120934228
$ws.Range("F17").Value = 0.9099865788991316
$ws.Range("G17").Value = 1.070171110887364
$ws.Range("H17").Value = 1.070171110887364
$ws.Range("I17").Value = 1.070171110887364
$ws.Range("J17").Value = 0.9801669931299923
$ws.Range("K17").Value = 1.032899157344269
$ws.Range("L17").Value = 1.025513505481809
$ws.Range("M17").Value = 0.9404269441365066
$ws.Range("N17").Value = 1.070171110887364
$ws.Range("O17").Value = 1.017051120934228
$ws.Range("P17").Value = 0.9635188499166798
$ws.Range("Q17").Value = 0.9986090570321102
$ws.Range("R17").Value = 0.9990696035735747
$ws.Range("S17").Value = 0.9690682309877839
$ws.Range("T17").Value = 0.9990696035735747
$ws.Range("U17").Value = 0.994343950962679
$ws.Range("V17").Value = 1.009509382947616
$ws.Range("W17").Value = 0.9943296626677578

# Row 18: Holden10
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Holden10"
$ws.Range("C18").Value = 0.9438331917532786
$ws.Range("D18").Value = 0.9917529270340689
$ws.Range("E18").Value = 1.018902682949478
$ws.Range("F18").Value = 0.9438331917532786
$ws.Range("G18").Value = 1.028962609637432
$ws.Range("H18").Value = 1.028962609637432
$ws.Range("I18").Value = 1.028962609637432
$ws.Range("J18").Value = 0.9881675675014191
$ws.Range("K18").Value = 1.021440756352737
$ws.Range("L18").Value = 1.011411455220193
$ws.Range("M18").Value = 0.9642337869449227
$ws.Range("N18").Value = 1.028962609637432
$ws.Range("O18").Value = 1.018902682949478
$ws.Range("P18").Value = 0.9813679373513784
$ws.Range("Q18").Value = 1.003535125225449
$ws.Range("R18").Value = 0.9972328281133963
$ws.Range("S18").Value = 0.9836344807347253
$ws.Range("T18").Value = 0.9972328281133963
$ws.Range("U18").Value = 0.9949665129604019
$ws.Range("V18").Value = 1.001765732295808
$ws.Range("W18").Value = 0.9960881221741912

# Row 19: Holden15
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Holden15"
$ws.Range("C19").Value = 0.9370261360452106
$ws.Range("D19").Value = 0.9920746955814405
$ws.Range("E19").Value = 1.023436297962328
$ws.Range("F19").Value = 0.9370261360452106
$ws.Range("G19").Value = 1.028753296839594
$ws.Range("H19").Value = 1.028753296839594
$ws.Range("I19").Value = 1.028753296839594
$ws.Range("J19").Value = 0.9867750525009706
$ws.Range("K19").Value = 1.024350514243677
$ws.Range("L19").Value = 1.011589993293249
$ws.Range("M19").Value = 0.9601931483567202
$ws.Range("N19").Value = 1.028753296839594
$ws.Range("O19").Value = 1.023436297962328
$ws.Range("P19").Value = 0.9802312170037691
$ws.Range("Q19").Value = 1.005105675231649
$ws.Range("R19").Value = 0.9964052436157109
$ws.Range("S19").Value = 0.9824124955028363
$ws.Range("T19").Value = 0.9964052436157109
$ws.Range("U19").Value = 0.9939976958370258
$ws.Range("V19").Value = 1.000948816037539
$ws.Range("W19").Value = 0.9955248918528987

# --- Rows 20-23: HexGrid rows (moved down from old 16-19; A column needs value + style copy) ---
# Row 20: HexGrid-90degTilt2.5degRes
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("C20").Value = 1.000014265262034
$ws.Range("D20").Value = 0.9999907513428428
$ws.Range("E20").Value = 1.000047888214149
$ws.Range("F20").Value = 1.000014265262034
$ws.Range("G20").Value = 1.000012951111389
$ws.Range("H20").Value = 1.000012951111389
$ws.Range("I20").Value = 1.000012951111389
$ws.Range("J20").Value = 0.9999664449995812
$ws.Range("K20").Value = 1.000030110702129
$ws.Range("L20").Value = 0.9999826819109803
$ws.Range("M20").Value = 0.9999800044695841
$ws.Range("N20").Value = 1.000012951111389
$ws.Range("O20").Value = 1.000047888214149
$ws.Range("P20").Value = 1.000031076738091
$ws.Range("Q20").Value = 1.000007166606865
$ws.Range("R20").Value = 1.000025034862524
$ws.Range("S20").Value = 1.000009532825255
$ws.Range("T20").Value = 1.000025034862524
$ws.Range("U20").Value = 1.000010387396788
$ws.Range("V20").Value = 1.000010900139709
$ws.Range("W20").Value = 1.000003137251586

# Row 21: HexGrid-90degTilt5degRes
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C21").Value = 0.9993719398151365
$ws.Range("D21").Value = 0.999947400678812
$ws.Range("E21").Value = 1.000191244423077
$ws.Range("F21").Value = 0.9993719398151365
$ws.Range("G21").Value = 1.000222050567939
$ws.Range("H21").Value = 1.000222050567939
$ws.Range("I21").Value = 1.000222050567939
$ws.Range("J21").Value = 0.99991443408531
$ws.Range("K21").Value = 1.000202620248907
$ws.Range("L21").Value = 1.000122387592581
$ws.Range("M21").Value = 0.9996416329785256
$ws.Range("N21").Value = 1.000222050567939
$ws.Range("O21").Value = 1.000191244423077
$ws.Range("P21").Value = 0.9997815921191067
$ws.Range("Q21").Value = 1.000052839254193
$ws.Range("R21").Value = 0.9999284116020507
$ws.Range("S21").Value = 0.9998258727745077
$ws.Range("T21").Value = 0.9999284116020505
$ws.Range("U21").Value = 0.9999249172228654
$ws.Range("V21").Value = 0.9999843438918801
$ws.Range("W21").Value = 0.9999517137987859

# Row 22: HexGrid-90degTilt10degRes
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("C22").Value = 0.998506415786273
$ws.Range("D22").Value = 0.9994819292093284
$ws.Range("E22").Value = 1.000221322261168
$ws.Range("F22").Value = 0.998506415786273
$ws.Range("G22").Value = 1.001586305192045
$ws.Range("H22").Value = 1.001586305192045
$ws.Range("I22").Value = 1.001586305192045
$ws.Range("J22").Value = 0.9995583097899157
$ws.Range("K22").Value = 1.000619838833334
$ws.Range("L22").Value = 1.000495921767031
$ws.Range("M22").Value = 0.9988977590113023
$ws.Range("N22").Value = 1.001586305192045
$ws.Range("O22").Value = 1.000221322261168
$ws.Range("P22").Value = 0.9993638690237203
$ws.Range("Q22").Value = 0.9998898160255416
$ws.Range("R22").Value = 1.000104681079828
$ws.Range("S22").Value = 0.9994286826124522
$ws.Range("T22").Value = 1.000104681079828
$ws.Range("U22").Value = 0.9999680882573503
$ws.Range("V22").Value = 1.000291731644289
$ws.Range("W22").Value = 0.9999209752312996

# Row 23: HexGrid-90degTilt15degRes
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"
$ws.Range("C23").Value = 0.9970591658172933
$ws.Range("D23").Value = 0.9985336736135125
$ws.Range("E23").Value = 1.000196223344717
$ws.Range("F23").Value = 0.9970591658172933
$ws.Range("G23").Value = 1.004322445273111
$ws.Range("H23").Value = 1.004322445273111
$ws.Range("I23").Value = 1.004322445273111
$ws.Range("J23").Value = 0.9988464451549904
$ws.Range("K23").Value = 1.001390740019026
$ws.Range("L23").Value = 1.001200027352157
$ws.Range("M23").Value = 0.9975362002347745
$ws.Range("N23").Value = 1.004322445273111
$ws.Range("O23").Value = 1.000196223344717
$ws.Range("P23").Value = 0.9986276945810049
$ws.Range("Q23").Value = 0.9995213342498535
$ws.Range("R23").Value = 1.000525944811707
$ws.Range("S23").Value = 0.9987006114390001
$ws.Range("T23").Value = 1.000525944811707
$ws.Range("U23").Value = 1.000106069897528
$ws.Range("V23").Value = 1.000949344972645
$ws.Range("W23").Value = 0.9998856151011977

# --- Apply column-A header style (bold, centered, bordered) to the new row labels (A16:A23) ---
$ws.Range("A2").Copy()
$ws.Range("A16:A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
